$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'91.284.75"
$ws.Range("E2").Value = "  +3.79%  "

# Row 3
$ws.Range("D3").Value = "'3.093.20"
$ws.Range("E3").Value = "  -0.50%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'218.12"
$ws.Range("E5").Value = "  +1.75%  "

# Row 6
$ws.Range("D6").Value = "'617.12"
$ws.Range("E6").Value = "  -2.90%  "

# Row 7
$ws.Range("D7").Value = "'0.372"
$ws.Range("E7").Value = "  -4.59%  "

# Row 8
$ws.Range("D8").Value = "'0.877"
$ws.Range("E8").Value = "  +11.22%  "

# Row 9
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("D10").Value = "'3.086.28"
$ws.Range("E10").Value = "  -0.63%  "

# Row 11
$ws.Range("D11").Value = "'0.664"
$ws.Range("E11").Value = "  +17.61%  "

# Row 12
$ws.Range("E12").Value = "  +5.75%  "

# Row 13
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("E13").Value = "  -0.14%  "

# Row 14
$ws.Range("D14").Value = "'90.892.98"
$ws.Range("E14").Value = "  +3.66%  "

# Row 15
$ws.Range("D15").Value = "'5.37"
$ws.Range("E15").Value = "  -0.28%  "

# Row 16
$ws.Range("D16").Value = "'32.92"
$ws.Range("E16").Value = "  +3.05%  "

# Row 17
$ws.Range("D17").Value = "'3.652.60"
$ws.Range("E17").Value = "  -0.66%  "

# Row 18
$ws.Range("D18").Value = "'3.080.74"
$ws.Range("E18").Value = "  -0.93%  "

# Row 19
$ws.Range("E19").Value = "  +4.29%  "

# Row 20
$ws.Range("D20").Value = "'0.0000219"
$ws.Range("E20").Value = "  -1.39%  "

# Row 21
$ws.Range("D21").Value = "'13.79"
$ws.Range("E21").Value = "  +4.31%  "

# Row 22
$ws.Range("D22").Value = "'434.40"
$ws.Range("E22").Value = "  +2.83%  "

# Row 23
$ws.Range("E23").Value = "  +0.39%  "

# Row 24
$ws.Range("D24").Value = "'5.11"

# Row 25
$ws.Range("D25").Value = "'5.64"
$ws.Range("E25").Value = "  +3.94%  "

# Row 26
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'83.82"
$ws.Range("E26").Value = "  +1.91%  "

# Row 27
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "'11.81"
$ws.Range("E27").Value = "  +2.89%  "

# Row 28
$ws.Range("D28").Value = "'3.242.01"
$ws.Range("E28").Value = "  -1.27%  "

# Row 30
$ws.Range("E30").Value = "  +8.04%  "

# Row 31
$ws.Range("D31").Value = "'0.166"
$ws.Range("E31").Value = "  +7.46%  "

# Row 32
$ws.Range("D32").Value = "'8.59"
$ws.Range("E32").Value = "  +5.11%  "

# Row 33
$ws.Range("D33").Value = "'3.84"
$ws.Range("E33").Value = "  -5.17%  "

# Row 34
$ws.Range("D34").Value = "'514.27"
$ws.Range("E34").Value = "  +2.55%  "

# Row 35
$ws.Range("D35").Value = "'6.89"
$ws.Range("E35").Value = "  +1.39%  "

# Row 36
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'1.27"
$ws.Range("E36").Value = "  -0.18%  "

# Row 37
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").Value = "'1.84"
$ws.Range("E37").Value = "  +0.39%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.137"
$ws.Range("E38").Value = "  -8.05%  "

# Row 39
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'23.01"
$ws.Range("E39").Value = "  +3.80%  "

# Row 40
$ws.Range("D40").Value = "'22.32"
$ws.Range("E40").Value = "  +0.60%  "

# Row 41
$ws.Range("E41").Value = "  +0.12%  "

# Row 43
$ws.Range("E43").Value = "  +4.32%  "

# Row 44
$ws.Range("E44").Value = "  +0.68%  "

# Row 45
$ws.Range("D45").Value = "'1.87"
$ws.Range("E45").Value = "  +1.53%  "

# Row 46
$ws.Range("D46").Value = "'0.0715"
$ws.Range("E46").Value = "  +11.01%  "

# Row 47
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'43.73"
$ws.Range("E47").Value = "  +0.22%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'141.08"
$ws.Range("E48").Value = "  -3.46%  "

# Row 49
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "'0.000263"
$ws.Range("E49").Value = "  +11.57%  "

# Row 50
$ws.Range("D50").Value = "'4.22"
$ws.Range("E50").Value = "  +7.17%  "

# Row 51
$ws.Range("D51").Value = "'165.01"
$ws.Range("E51").Value = "  +1.81%  "
